$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2021-12-17"

# Update header label in B1
$ws.Range("B1").Value = "December 2021 (through December 17)"

# Update individual cell values per the new data snapshot
$ws.Range("AX3").Value = 2
$ws.Range("BV4").Value = 2
$ws.Range("AX5").Value = 2
$ws.Range("BJ6").Value = 5
$ws.Range("AL7").Value = 8
$ws.Range("BJ8").Value = 3
$ws.Range("BV11").Value = 2
$ws.Range("AX16").Value = 3
$ws.Range("AX18").Value = 2
$ws.Range("AL22").Value = 2
$ws.Range("AX22").Value = 3
$ws.Range("N23").Value = 1
$ws.Range("AX24").Value = 4
$ws.Range("BJ24").Value = 4
$ws.Range("BJ30").Value = 2
$ws.Range("B37").Value = 3
$ws.Range("AX64").Value = 1
$ws.Range("AL66").Value = 2
$ws.Range("AX67").Value = 1
$ws.Range("AL78").Value = 1
$ws.Range("N82").Value = 4
$ws.Range("AX84").Value = 2
